$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("H94").Value = 636
$ws.Range("I94").Value = 379.63635
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 379.63635
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = 71.36365000000001
$ws.Range("N94").Value = -2102
$ws.Range("H137").Value = 2920.1428
$ws.Range("I137").Value = 3133.077
$ws.Range("J137").Value = 2574.125
$ws.Range("K137").Value = 9399.231
$ws.Range("L137").Value = 7722.375
$ws.Range("M137").Value = -6849.231
$ws.Range("N137").Value = -12822.375
$ws.Range("M64").Value = -4752
$ws.Range("M67").Value = -4142

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13178.511
$ws.Range("I32").Value = 8046.057
$ws.Range("K32").Value = 8046.057
$ws.Range("M32").Value = -7759.057
$ws.Range("H45").Value = 2113.52
$ws.Range("J45").Value = 2929.9
$ws.Range("L45").Value = 2929.9
$ws.Range("N45").Value = -3683.9
$ws.Range("H61").Value = 2500.3677
$ws.Range("J61").Value = 3945.7273
$ws.Range("L61").Value = 3945.7273
$ws.Range("N61").Value = -4369.7273
$ws.Range("H74").Value = 4542.2563
$ws.Range("I74").Value = 2108.4517
$ws.Range("K74").Value = 2108.4517
$ws.Range("M74").Value = -1234.4517
$ws.Range("H77").Value = 4542.2563
$ws.Range("I77").Value = 2108.4517
$ws.Range("K77").Value = 10542.2585
$ws.Range("M77").Value = -6174.2585
$ws.Range("H132").Value = 6844.6
$ws.Range("I132").Value = 5764.0835
$ws.Range("K132").Value = 17292.2505
$ws.Range("M132").Value = -14762.2505
$ws.Range("H136").Value = 2500.3677
$ws.Range("J136").Value = 3945.7273
$ws.Range("L136").Value = 11837.1819
$ws.Range("N136").Value = -16937.1819

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7165.3184
$ws.Range("I86").Value = 5616.8
$ws.Range("J86").Value = 10483.571
$ws.Range("K86").Value = 5616.8
$ws.Range("L86").Value = 10483.571
$ws.Range("M86").Value = -4493.8
$ws.Range("N86").Value = -12729.571
$ws.Range("H89").Value = 7165.3184
$ws.Range("I89").Value = 5616.8
$ws.Range("J89").Value = 10483.571
$ws.Range("K89").Value = 28084
$ws.Range("L89").Value = 52417.855
$ws.Range("M89").Value = -22468
$ws.Range("N89").Value = -63649.855
$ws.Range("H92").Value = 54200
$ws.Range("J92").Value = 54200
$ws.Range("L92").Value = 54200
$ws.Range("N92").Value = -59192
$ws.Range("H107").Value = 2540.8
$ws.Range("J107").Value = 3634.6667
$ws.Range("L107").Value = 3634.6667
$ws.Range("N107").Value = -7474.6667

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10764.05
$ws.Range("I31").Value = 3244.9
$ws.Range("K31").Value = 3244.9
$ws.Range("M31").Value = -2949.9
$ws.Range("H34").Value = 10764.05
$ws.Range("I34").Value = 3244.9
$ws.Range("K34").Value = 3244.9
$ws.Range("M34").Value = -3042.9
$ws.Range("H58").Value = 2027.9697
$ws.Range("I58").Value = 1338.5
$ws.Range("K58").Value = 1338.5
$ws.Range("M58").Value = -1135.5
$ws.Range("H105").Value = 1637.5555
$ws.Range("J105").Value = 564
$ws.Range("L105").Value = 564
$ws.Range("N105").Value = -4058
$ws.Range("H134").Value = 2248.3
$ws.Range("I134").Value = 2038.9688
$ws.Range("K134").Value = 6116.9064
$ws.Range("M134").Value = -3581.9064
$ws.Range("H136").Value = 2027.9697
$ws.Range("I136").Value = 1338.5
$ws.Range("K136").Value = 4015.5
$ws.Range("M136").Value = -1465.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1642.6666
$ws.Range("I18").Value = 499
$ws.Range("J18").Value = 2214.5
$ws.Range("K18").Value = 1497
$ws.Range("L18").Value = 6643.5
$ws.Range("M18").Value = -1328
$ws.Range("N18").Value = -6981.5
$ws.Range("H49").Value = 1409.4
$ws.Range("J49").Value = 1049.5
$ws.Range("L49").Value = 3148.5
$ws.Range("N49").Value = -3460.5
$ws.Range("H51").Value = 1500
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H68").Value = 1765.7778
$ws.Range("I68").Value = 1318.4
$ws.Range("J68").Value = 2325
$ws.Range("K68").Value = 3955.2
$ws.Range("L68").Value = 6975
$ws.Range("M68").Value = -3144.2
$ws.Range("N68").Value = -8597
$ws.Range("H71").Value = 1765.7778
$ws.Range("I71").Value = 1318.4
$ws.Range("J71").Value = 2325
$ws.Range("K71").Value = 11865.6
$ws.Range("L71").Value = 20925
$ws.Range("M71").Value = -7809.6
$ws.Range("N71").Value = -29037
$ws.Range("H86").Value = 7060.4443
$ws.Range("J86").Value = 2191.1667
$ws.Range("L86").Value = 6573.500100000001
$ws.Range("N86").Value = -8945.500100000001
$ws.Range("H89").Value = 7060.4443
$ws.Range("J89").Value = 2191.1667
$ws.Range("L89").Value = 19720.5003
$ws.Range("N89").Value = -31576.5003
$ws.Range("H92").Value = 4349.6665
$ws.Range("I92").Value = 4000
$ws.Range("J92").Value = 4524.5
$ws.Range("K92").Value = 12000
$ws.Range("L92").Value = 13573.5
$ws.Range("M92").Value = -10752
$ws.Range("N92").Value = -16069.5
$ws.Range("H97").Value = 2317.724
$ws.Range("I97").Value = 2920.65
$ws.Range("J97").Value = 977.8889
$ws.Range("K97").Value = 8761.950000000001
$ws.Range("L97").Value = 2933.6667
$ws.Range("M97").Value = -8265.950000000001
$ws.Range("N97").Value = -3925.6667
$ws.Range("H104").Value = 7700
$ws.Range("J104").Value = 7700
$ws.Range("L104").Value = 23100
$ws.Range("N104").Value = -28342
$ws.Range("H107").Value = 237.66667
$ws.Range("I107").Value = 75
$ws.Range("J107").Value = 296.81818
$ws.Range("K107").Value = 225
$ws.Range("L107").Value = 890.45454
$ws.Range("M107").Value = 1695
$ws.Range("N107").Value = -4730.45454

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6494.543
$ws.Range("I70").Value = 5547.4
$ws.Range("K70").Value = 5547.4
$ws.Range("M70").Value = -5277.4
$ws.Range("H73").Value = 6494.543
$ws.Range("I73").Value = 5547.4
$ws.Range("K73").Value = 5547.4
$ws.Range("M73").Value = -4611.4
$ws.Range("H80").Value = 2710.0588
$ws.Range("I80").Value = 2703.8572
$ws.Range("K80").Value = 2703.8572
$ws.Range("M80").Value = -1705.8572
$ws.Range("H83").Value = 2710.0588
$ws.Range("I83").Value = 2703.8572
$ws.Range("K83").Value = 13519.286
$ws.Range("M83").Value = -8527.286
$ws.Range("H126").Value = 5768.5
$ws.Range("I126").Value = 6756.778
$ws.Range("J126").Value = 4497.857
$ws.Range("K126").Value = 20270.334
$ws.Range("L126").Value = 13493.571
$ws.Range("M126").Value = -17800.334
$ws.Range("N126").Value = -18433.571
$ws.Range("H132").Value = 5286.913
$ws.Range("I132").Value = 4905.905
$ws.Range("K132").Value = 14717.715
$ws.Range("M132").Value = -12187.715

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4831.875
$ws.Range("I7").Value = 4831.875
$ws.Range("K7").Value = 4831.875
$ws.Range("M7").Value = -4719.875
$ws.Range("H40").Value = 5361.375
$ws.Range("I40").Value = 4598.5835
$ws.Range("K40").Value = 4598.5835
$ws.Range("M40").Value = -4462.5835
$ws.Range("H68").Value = 6084.478
$ws.Range("I68").Value = 5457.1665
$ws.Range("J68").Value = 6768.8184
$ws.Range("K68").Value = 5457.1665
$ws.Range("L68").Value = 6768.8184
$ws.Range("M68").Value = -4708.1665
$ws.Range("N68").Value = -8266.8184
$ws.Range("H71").Value = 6084.478
$ws.Range("I71").Value = 5457.1665
$ws.Range("J71").Value = 6768.8184
$ws.Range("K71").Value = 27285.8325
$ws.Range("L71").Value = 33844.092
$ws.Range("M71").Value = -23541.8325
$ws.Range("N71").Value = -41332.092
$ws.Range("H122").Value = 5050.5713
$ws.Range("I122").Value = 3159.4285
$ws.Range("J122").Value = 6941.7144
$ws.Range("K122").Value = 9478.2855
$ws.Range("L122").Value = 20825.1432
$ws.Range("M122").Value = -7028.2855
$ws.Range("N122").Value = -25725.1432
$ws.Range("H126").Value = 4831.875
$ws.Range("I126").Value = 4831.875
$ws.Range("K126").Value = 14495.625
$ws.Range("M126").Value = -12025.625
$ws.Range("H132").Value = 7277.4287
$ws.Range("I132").Value = 6471.75
$ws.Range("J132").Value = 7599.7
$ws.Range("K132").Value = 19415.25
$ws.Range("L132").Value = 22799.1
$ws.Range("M132").Value = -16885.25
$ws.Range("N132").Value = -27859.1
$ws.Range("H136").Value = 7070.794
$ws.Range("I136").Value = 2635.5356
$ws.Range("K136").Value = 7906.6068
$ws.Range("M136").Value = -5356.6068

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6789.5713
$ws.Range("I132").Value = 6584
$ws.Range("K132").Value = 19752
$ws.Range("M132").Value = -17222
$ws.Range("H136").Value = 1732.0492
$ws.Range("I136").Value = 1482.5209
$ws.Range("K136").Value = 4447.5627
$ws.Range("M136").Value = -1897.5627
